# Auto-generated cell updates derived from the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 396.5
$ws.Range("I55").Value = 396.42856
$ws.Range("J55").Value = 396.66666
$ws.Range("K55").Value = 396.42856
$ws.Range("L55").Value = 396.66666
$ws.Range("M55").Value = -182.42856
$ws.Range("N55").Value = -824.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2783.9736
$ws.Range("I138").Value = 2217.4736
$ws.Range("J138").Value = 3350.4736
$ws.Range("K138").Value = 6652.4208
$ws.Range("L138").Value = 10051.4208
$ws.Range("M138").Value = -1512.4208
$ws.Range("N138").Value = -20331.4208

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2546.5686
$ws.Range("I141").Value = 1116.375
$ws.Range("J141").Value = 4955.316
$ws.Range("K141").Value = 3349.125
$ws.Range("L141").Value = 14865.948
$ws.Range("M141").Value = 1830.875
$ws.Range("N141").Value = -25225.948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8517.724
$ws.Range("I32").Value = 5097.9077
$ws.Range("J32").Value = 28725.727
$ws.Range("K32").Value = 5097.9077
$ws.Range("L32").Value = 28725.727
$ws.Range("M32").Value = -4810.9077
$ws.Range("N32").Value = -29299.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3321.5715
$ws.Range("I45").Value = 2212.2
$ws.Range("J45").Value = 3937.889
$ws.Range("K45").Value = 2212.2
$ws.Range("L45").Value = 3937.889
$ws.Range("M45").Value = -1835.2
$ws.Range("N45").Value = -4691.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 38500
$ws.Range("J54").Value = 38500
$ws.Range("L54").Value = 38500
$ws.Range("N54").Value = -40038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 137499.25
$ws.Range("I80").Value = 100000
$ws.Range("J80").Value = 149999
$ws.Range("K80").Value = 100000
$ws.Range("L80").Value = 149999
$ws.Range("M80").Value = -99002
$ws.Range("N80").Value = -151995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 137499.25
$ws.Range("I83").Value = 100000
$ws.Range("J83").Value = 149999
$ws.Range("K83").Value = 300000
$ws.Range("L83").Value = 449997
$ws.Range("M83").Value = -295008
$ws.Range("N83").Value = -459981

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2039.1
$ws.Range("I107").Value = 2067
$ws.Range("J107").Value = 2027.1428
$ws.Range("K107").Value = 2067
$ws.Range("L107").Value = 2027.1428
$ws.Range("M107").Value = -147
$ws.Range("N107").Value = -5867.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2707.4028
$ws.Range("I134").Value = 1969.194
$ws.Range("J134").Value = 12599.4
$ws.Range("K134").Value = 5907.582
$ws.Range("L134").Value = 37798.2
$ws.Range("M134").Value = -3372.582
$ws.Range("N134").Value = -42868.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 90000
$ws.Range("I139").Value = 60000
$ws.Range("J139").Value = 120000
$ws.Range("K139").Value = 60000
$ws.Range("L139").Value = 120000
$ws.Range("M139").Value = -54860
$ws.Range("N139").Value = -130280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 881.7143
$ws.Range("I16").Value = 887.8333
$ws.Range("J16").Value = 845
$ws.Range("K16").Value = 887.8333
$ws.Range("L16").Value = 845
$ws.Range("M16").Value = -600.8333
$ws.Range("N16").Value = -1419

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1502.9149
$ws.Range("I58").Value = 913
$ws.Range("J58").Value = 2893.4285
$ws.Range("K58").Value = 913
$ws.Range("L58").Value = 2893.4285
$ws.Range("M58").Value = -710
$ws.Range("N58").Value = -3299.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1578.7736
$ws.Range("I99").Value = 1619.5349
$ws.Range("J99").Value = 1403.5
$ws.Range("K99").Value = 1619.5349
$ws.Range("L99").Value = 1403.5
$ws.Range("M99").Value = -121.5349000000001
$ws.Range("N99").Value = -4399.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1732.8889
$ws.Range("I105").Value = 2029.5714
$ws.Range("K105").Value = 2029.5714
$ws.Range("M105").Value = -282.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1696.3549
$ws.Range("I107").Value = 1414.8966
$ws.Range("J107").Value = 5777.5
$ws.Range("K107").Value = 1414.8966
$ws.Range("L107").Value = 5777.5
$ws.Range("M107").Value = 505.1034
$ws.Range("N107").Value = -9617.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 881.7143
$ws.Range("I113").Value = 887.8333
$ws.Range("J113").Value = 845
$ws.Range("K113").Value = 887.8333
$ws.Range("L113").Value = 845
$ws.Range("M113").Value = 1282.1667
$ws.Range("N113").Value = -5185

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1578.7736
$ws.Range("I126").Value = 1619.5349
$ws.Range("J126").Value = 1403.5
$ws.Range("K126").Value = 4858.6047
$ws.Range("L126").Value = 4210.5
$ws.Range("M126").Value = -2388.6047
$ws.Range("N126").Value = -9150.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2540.0857
$ws.Range("I132").Value = 1718.875
$ws.Range("K132").Value = 5156.625
$ws.Range("M132").Value = -2626.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1502.9149
$ws.Range("I136").Value = 913
$ws.Range("J136").Value = 2893.4285
$ws.Range("K136").Value = 2739
$ws.Range("L136").Value = 8680.2855
$ws.Range("M136").Value = -189
$ws.Range("N136").Value = -13780.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 417.81818
$ws.Range("I2").Value = 596.5714
$ws.Range("K2").Value = 3579.4284
$ws.Range("M2").Value = -3466.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3716258.5
$ws.Range("I4").Value = 3816694.2
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 11450082.6
$ws.Range("L4").Value = 420
$ws.Range("M4").Value = -11449970.6
$ws.Range("N4").Value = -644

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5773.4707
$ws.Range("I34").Value = 580
$ws.Range("K34").Value = 1740
$ws.Range("M34").Value = -1656

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1874.4166
$ws.Range("J46").Value = 2456.2856
$ws.Range("L46").Value = 7368.8568
$ws.Range("N46").Value = -7550.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2249
$ws.Range("I54").Value = 1499
$ws.Range("J54").Value = 2999
$ws.Range("K54").Value = 4497
$ws.Range("L54").Value = 8997
$ws.Range("M54").Value = -3938
$ws.Range("N54").Value = -10115

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6493.5
$ws.Range("J75").Value = 6493.5
$ws.Range("L75").Value = 19480.5
$ws.Range("N75").Value = -21476.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 6493.5
$ws.Range("J78").Value = 6493.5
$ws.Range("L78").Value = 58441.5
$ws.Range("N78").Value = -68425.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 940.5
$ws.Range("I109").Value = 940.5
$ws.Range("K109").Value = 2821.5
$ws.Range("M109").Value = -1781.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1272
$ws.Range("I114").Value = 500
$ws.Range("J114").Value = 1786.6666
$ws.Range("K114").Value = 1500
$ws.Range("L114").Value = 5359.9998
$ws.Range("M114").Value = 1754
$ws.Range("N114").Value = -11867.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4247.684
$ws.Range("I137").Value = 2882.1667
$ws.Range("J137").Value = 4503.7188
$ws.Range("K137").Value = 8646.500100000001
$ws.Range("L137").Value = 13511.1564
$ws.Range("M137").Value = -3546.500100000001
$ws.Range("N137").Value = -23711.1564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4071.2222
$ws.Range("I139").Value = 4090.3333
$ws.Range("J139").Value = 4033
$ws.Range("K139").Value = 12270.9999
$ws.Range("L139").Value = 12099
$ws.Range("M139").Value = -7130.999899999999
$ws.Range("N139").Value = -22379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 571.5454999999999
$ws.Range("J107").Value = 998.75
$ws.Range("L107").Value = 998.75
$ws.Range("N107").Value = -4838.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3686.1667
$ws.Range("I113").Value = 3574
$ws.Range("J113").Value = 3910.5
$ws.Range("K113").Value = 3574
$ws.Range("L113").Value = 3910.5
$ws.Range("M113").Value = -1404
$ws.Range("N113").Value = -8250.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3131.8572
$ws.Range("I132").Value = 2747.432
$ws.Range("K132").Value = 8242.295999999998
$ws.Range("M132").Value = -5712.295999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3865.3572
$ws.Range("I132").Value = 2729.5312
$ws.Range("K132").Value = 8188.5936
$ws.Range("M132").Value = -5658.5936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7181.3257
$ws.Range("I136").Value = 3019.926
$ws.Range("J136").Value = 9085.695
$ws.Range("K136").Value = 9059.778
$ws.Range("L136").Value = 27257.085
$ws.Range("M136").Value = -6509.778
$ws.Range("N136").Value = -32357.085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 58333
$ws.Range("I51").Value = 31999
$ws.Range("J51").Value = 71500
$ws.Range("K51").Value = 31999
$ws.Range("L51").Value = 71500
$ws.Range("M51").Value = -31489
$ws.Range("N51").Value = -72520

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1393.1052
$ws.Range("I107").Value = 1345.5
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 4036.5
$ws.Range("L107").Value = 6750
$ws.Range("M107").Value = -2116.5
$ws.Range("N107").Value = -10590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 841.1429000000001
$ws.Range("I113").Value = 841.1429000000001
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2523.4287
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -353.4287000000004
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1518.3855
$ws.Range("I136").Value = 1407.8197
$ws.Range("J136").Value = 1824.9546
$ws.Range("K136").Value = 4223.4591
$ws.Range("L136").Value = 5474.8638
$ws.Range("M136").Value = -1673.4591
$ws.Range("N136").Value = -10574.8638
